$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chessboard coordinate labels: "wit:" (white squares) rows 16-23, "grijs:" (gray squares) rows 25-32
$ws.Range("F16").Value = "wit:"
$ws.Range("G16").Value = "r=0,c=0"
$ws.Range("H16").Value = "r=0,c=2"
$ws.Range("I16").Value = "r=0,c=4"
$ws.Range("J16").Value = "r=0,c=6"

$ws.Range("G17").Value = "r=1,c=1"
$ws.Range("H17").Value = "r=1,c=3"
$ws.Range("I17").Value = "r=1,c=5"
$ws.Range("J17").Value = "r=1,c=7"

$ws.Range("G18").Value = "r=2,c=0"
$ws.Range("H18").Value = "r=2,c=2"
$ws.Range("I18").Value = "r=2,c=4"
$ws.Range("J18").Value = "r=2,c=6"

$ws.Range("G19").Value = "r=3,c=1"
$ws.Range("H19").Value = "r=3,c=3"
$ws.Range("I19").Value = "r=3,c=5"
$ws.Range("J19").Value = "r=3,c=7"

$ws.Range("G20").Value = "r=4,c=0"
$ws.Range("H20").Value = "r=4,c=2"
$ws.Range("I20").Value = "r=4,c=4"
$ws.Range("J20").Value = "r=4,c=6"

$ws.Range("G21").Value = "r=5,c=1"
$ws.Range("H21").Value = "r=5,c=3"
$ws.Range("I21").Value = "r=5,c=5"
$ws.Range("J21").Value = "r=5,c=7"

$ws.Range("G22").Value = "r=6,c=0"
$ws.Range("H22").Value = "r=6,c=2"
$ws.Range("I22").Value = "r=6,c=4"
$ws.Range("J22").Value = "r=6,c=6"

$ws.Range("G23").Value = "r=7,c=1"
$ws.Range("H23").Value = "r=7,c=3"
$ws.Range("I23").Value = "r=7,c=5"
$ws.Range("J23").Value = "r=7,c=7"

$ws.Range("F25").Value = "grijs:"
$ws.Range("G25").Value = "r=0,c=1"
$ws.Range("H25").Value = "r=0,c=3"
$ws.Range("I25").Value = "r=0,c=5"
$ws.Range("J25").Value = "r=0,c=7"

$ws.Range("G26").Value = "r=1,c=0"
$ws.Range("H26").Value = "r=1,c=2"
$ws.Range("I26").Value = "r=1,c=4"
$ws.Range("J26").Value = "r=1,c=6"

$ws.Range("G27").Value = "r=2,c=1"
$ws.Range("H27").Value = "r=2,c=3"
$ws.Range("I27").Value = "r=2,c=5"
$ws.Range("J27").Value = "r=2,c=7"

$ws.Range("G28").Style = "Normal"
$ws.Range("G28").Value = "r=3,c=0"
$ws.Range("H28").Value = "r=3,c=2"
$ws.Range("I28").Value = "r=3,c=4"
$ws.Range("J28").Value = "r=3,c=6"

$ws.Range("G29").Style = "Normal"
$ws.Range("G29").Value = "r=4,c=1"
$ws.Range("H29").Value = "r=4,c=3"
$ws.Range("I29").Value = "r=4,c=5"
$ws.Range("J29").Value = "r=4,c=7"

$ws.Range("G30").Style = "Normal"
$ws.Range("G30").Value = "r=5,c=0"
$ws.Range("H30").Value = "r=5,c=2"
$ws.Range("I30").Value = "r=5,c=4"
$ws.Range("J30").Value = "r=5,c=6"

$ws.Range("G31").Style = "Normal"
$ws.Range("G31").Value = "r=6,c=1"
$ws.Range("H31").Value = "r=6,c=3"
$ws.Range("I31").Value = "r=6,c=5"
$ws.Range("J31").Value = "r=6,c=7"

$ws.Range("G32").Style = "Normal"
$ws.Range("G32").Value = "r=7,c=0"
$ws.Range("H32").Value = "r=7,c=2"
$ws.Range("I32").Value = "r=7,c=4"
$ws.Range("J32").Value = "r=7,c=6"

# Numeric sum grid rows 34-50
$ws.Range("F34").ClearFormats()
$ws.Range("F34").Value = "wit:"
$ws.Range("G34").Style = "Normal"
$ws.Range("G34").Value = 0
$ws.Range("H34").Style = "Normal"
$ws.Range("H34").Value = 2
$ws.Range("I34").Style = "Normal"
$ws.Range("I34").Value = 4
$ws.Range("J34").Style = "Normal"
$ws.Range("J34").Value = 6

$ws.Range("F35").Clear()
$ws.Range("G35").Style = "Normal"
$ws.Range("G35").Value = 2
$ws.Range("H35").Style = "Normal"
$ws.Range("H35").Value = 4
$ws.Range("I35").Style = "Normal"
$ws.Range("I35").Value = 6
$ws.Range("J35").Style = "Normal"
$ws.Range("J35").Value = 8

$ws.Range("F36").Clear()
$ws.Range("G36").Style = "Normal"
$ws.Range("G36").Value = 2
$ws.Range("H36").Style = "Normal"
$ws.Range("H36").Value = 4
$ws.Range("I36").Style = "Normal"
$ws.Range("I36").Value = 6
$ws.Range("J36").Style = "Normal"
$ws.Range("J36").Value = 8

$ws.Range("F37").Clear()
$ws.Range("G37").Style = "Normal"
$ws.Range("G37").Value = 4
$ws.Range("H37").Style = "Normal"
$ws.Range("H37").Value = 6
$ws.Range("I37").Style = "Normal"
$ws.Range("I37").Value = 8
$ws.Range("J37").Style = "Normal"
$ws.Range("J37").Value = 10

$ws.Range("F38").Clear()
$ws.Range("G38").Style = "Normal"
$ws.Range("G38").Value = 4
$ws.Range("H38").Style = "Normal"
$ws.Range("H38").Value = 6
$ws.Range("I38").Style = "Normal"
$ws.Range("I38").Value = 8
$ws.Range("J38").Style = "Normal"
$ws.Range("J38").Value = 10

$ws.Range("G39").Style = "Normal"
$ws.Range("G39").Value = 6
$ws.Range("H39").Style = "Normal"
$ws.Range("H39").Value = 8
$ws.Range("I39").Style = "Normal"
$ws.Range("I39").Value = 10
$ws.Range("J39").Style = "Normal"
$ws.Range("J39").Value = 12

$ws.Range("F40").Clear()
$ws.Range("G40").Style = "Normal"
$ws.Range("G40").Value = 6
$ws.Range("H40").Style = "Normal"
$ws.Range("H40").Value = 8
$ws.Range("I40").Style = "Normal"
$ws.Range("I40").Value = 10
$ws.Range("J40").Style = "Normal"
$ws.Range("J40").Value = 12

$ws.Range("F41").Clear()
$ws.Range("G41").Style = "Normal"
$ws.Range("G41").Value = 8
$ws.Range("H41").Style = "Normal"
$ws.Range("H41").Value = 10
$ws.Range("I41").Style = "Normal"
$ws.Range("I41").Value = 12
$ws.Range("J41").Style = "Normal"
$ws.Range("J41").Value = 14

$ws.Range("F42").Clear()
$ws.Range("G42").Clear()

$ws.Range("F43").ClearFormats()
$ws.Range("F43").Value = "grijs:"
$ws.Range("G43").Style = "Normal"
$ws.Range("G43").Value = 1
$ws.Range("H43").Style = "Normal"
$ws.Range("H43").Value = 3
$ws.Range("I43").Style = "Normal"
$ws.Range("I43").Value = 5
$ws.Range("J43").Style = "Normal"
$ws.Range("J43").Value = 7

$ws.Range("F44").Clear()
$ws.Range("G44").Style = "Normal"
$ws.Range("G44").Value = 1
$ws.Range("H44").Style = "Normal"
$ws.Range("H44").Value = 3
$ws.Range("I44").Style = "Normal"
$ws.Range("I44").Value = 5
$ws.Range("J44").Style = "Normal"
$ws.Range("J44").Value = 6

$ws.Range("G45").Style = "Normal"
$ws.Range("G45").Value = 3
$ws.Range("H45").Style = "Normal"
$ws.Range("H45").Value = 5
$ws.Range("I45").Style = "Normal"
$ws.Range("I45").Value = 7
$ws.Range("J45").Style = "Normal"
$ws.Range("J45").Value = 9

$ws.Range("G46").Style = "Normal"
$ws.Range("G46").Value = 3
$ws.Range("H46").Style = "Normal"
$ws.Range("H46").Value = 5
$ws.Range("I46").Style = "Normal"
$ws.Range("I46").Value = 7
$ws.Range("J46").Style = "Normal"
$ws.Range("J46").Value = 9

$ws.Range("G47").Style = "Normal"
$ws.Range("G47").Value = 5
$ws.Range("H47").Style = "Normal"
$ws.Range("H47").Value = 7
$ws.Range("I47").Style = "Normal"
$ws.Range("I47").Value = 9
$ws.Range("J47").Style = "Normal"
$ws.Range("J47").Value = 11

$ws.Range("G48").Style = "Normal"
$ws.Range("G48").Value = 5
$ws.Range("H48").Style = "Normal"
$ws.Range("H48").Value = 7
$ws.Range("I48").Style = "Normal"
$ws.Range("I48").Value = 9
$ws.Range("J48").Style = "Normal"
$ws.Range("J48").Value = 11

$ws.Range("G49").Style = "Normal"
$ws.Range("G49").Value = 7
$ws.Range("H49").Style = "Normal"
$ws.Range("H49").Value = 9
$ws.Range("I49").Style = "Normal"
$ws.Range("I49").Value = 11
$ws.Range("J49").Style = "Normal"
$ws.Range("J49").Value = 13

$ws.Range("G50").Style = "Normal"
$ws.Range("G50").Value = 7
$ws.Range("H50").Style = "Normal"
$ws.Range("H50").Value = 9
$ws.Range("I50").Style = "Normal"
$ws.Range("I50").Value = 11
$ws.Range("J50").Style = "Normal"
$ws.Range("J50").Value = 13

# Update selection
$ws.Range("K12").Select()